$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 760.6896400000001
$ws.Range("I33").Value = 609.2857
$ws.Range("J33").Value = 5000
$ws.Range("K33").Value = 609.2857
$ws.Range("L33").Value = 5000
$ws.Range("M33").Value = -380.2857
$ws.Range("N33").Value = -5458

$ws.Range("H103").Value = 3242.8572
$ws.Range("I103").Value = 4106.6665
$ws.Range("J103").Value = 1083.3334
$ws.Range("K103").Value = 12319.9995
$ws.Range("L103").Value = 3250.0002
$ws.Range("M103").Value = -11733.9995
$ws.Range("N103").Value = -4422.0002

$ws.Range("H125").Value = 2444.353
$ws.Range("I125").Value = 1786.3334
$ws.Range("J125").Value = 3184.625
$ws.Range("K125").Value = 16077.0006
$ws.Range("L125").Value = 28661.625
$ws.Range("M125").Value = -13617.0006
$ws.Range("N125").Value = -33581.625

$ws.Range("H135").Value = 1000
$ws.Range("J135").Value = 1000
$ws.Range("L135").Value = 9000
$ws.Range("N135").Value = -14070

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4850.0557
$ws.Range("I32").Value = 2887.814
$ws.Range("J32").Value = 12520.637
$ws.Range("K32").Value = 2887.814
$ws.Range("L32").Value = 12520.637
$ws.Range("M32").Value = -2600.814
$ws.Range("N32").Value = -13094.637

$ws.Range("H61").Value = 1504.138
$ws.Range("I61").Value = 1306.4
$ws.Range("J61").Value = 2740
$ws.Range("K61").Value = 1306.4
$ws.Range("L61").Value = 2740
$ws.Range("M61").Value = -1094.4
$ws.Range("N61").Value = -3164

$ws.Range("H74").Value = 43231.27
$ws.Range("I74").Value = 72363.21000000001
$ws.Range("J74").Value = 9244
$ws.Range("K74").Value = 72363.21000000001
$ws.Range("L74").Value = 9244
$ws.Range("M74").Value = -71489.21000000001
$ws.Range("N74").Value = -10992

$ws.Range("H77").Value = 43231.27
$ws.Range("I77").Value = 72363.21000000001
$ws.Range("J77").Value = 9244
$ws.Range("K77").Value = 361816.05
$ws.Range("L77").Value = 46220
$ws.Range("M77").Value = -357448.05
$ws.Range("N77").Value = -54956

$ws.Range("H132").Value = 4572.857
$ws.Range("I132").Value = 4378
$ws.Range("J132").Value = 4832.6665
$ws.Range("K132").Value = 13134
$ws.Range("L132").Value = 14497.9995
$ws.Range("M132").Value = -10604
$ws.Range("N132").Value = -19557.9995

$ws.Range("H136").Value = 1504.138
$ws.Range("I136").Value = 1306.4
$ws.Range("J136").Value = 2740
$ws.Range("K136").Value = 3919.2
$ws.Range("L136").Value = 8220
$ws.Range("M136").Value = -1369.2
$ws.Range("N136").Value = -13320

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 757.875
$ws.Range("I94").Value = 702.36365
$ws.Range("J94").Value = 880
$ws.Range("K94").Value = 702.36365
$ws.Range("L94").Value = 880
$ws.Range("M94").Value = -251.36365
$ws.Range("N94").Value = -1782

$ws.Range("H99").Value = 2021.7
$ws.Range("I99").Value = 1869.5
$ws.Range("K99").Value = 1869.5
$ws.Range("M99").Value = -371.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1064.2253
$ws.Range("I31").Value = 922.8205
$ws.Range("J31").Value = 1236.5625
$ws.Range("K31").Value = 922.8205
$ws.Range("L31").Value = 1236.5625
$ws.Range("M31").Value = -627.8205
$ws.Range("N31").Value = -1826.5625

$ws.Range("H34").Value = 1064.2253
$ws.Range("I34").Value = 922.8205
$ws.Range("J34").Value = 1236.5625
$ws.Range("K34").Value = 922.8205
$ws.Range("L34").Value = 1236.5625
$ws.Range("M34").Value = -720.8205
$ws.Range("N34").Value = -1640.5625

$ws.Range("H99").Value = 5954306
$ws.Range("I99").Value = 2233.3333
$ws.Range("J99").Value = 23810524
$ws.Range("K99").Value = 2233.3333
$ws.Range("L99").Value = 23810524
$ws.Range("M99").Value = -735.3332999999998
$ws.Range("N99").Value = -23813520

$ws.Range("H107").Value = 931.0625
$ws.Range("I107").Value = 498.22223
$ws.Range("J107").Value = 1487.5714
$ws.Range("K107").Value = 498.22223
$ws.Range("L107").Value = 1487.5714
$ws.Range("M107").Value = 1421.77777
$ws.Range("N107").Value = -5327.5714

$ws.Range("H126").Value = 5954306
$ws.Range("I126").Value = 2233.3333
$ws.Range("J126").Value = 23810524
$ws.Range("K126").Value = 6699.999899999999
$ws.Range("L126").Value = 71431572
$ws.Range("M126").Value = -4229.999899999999
$ws.Range("N126").Value = -71436512

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 159.61539
$ws.Range("J7").Value = 177.63637
$ws.Range("L7").Value = 532.9091100000001
$ws.Range("N7").Value = -756.9091100000001

$ws.Range("H12").Value = 374.3684
$ws.Range("I12").Value = 428.625
$ws.Range("J12").Value = 334.9091
$ws.Range("K12").Value = 1285.875
$ws.Range("L12").Value = 1004.7273
$ws.Range("M12").Value = -1112.875
$ws.Range("N12").Value = -1350.7273

$ws.Range("H33").Value = 67.75
$ws.Range("J33").Value = 67.75
$ws.Range("L33").Value = 406.5
$ws.Range("N33").Value = -972.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5425.552
$ws.Range("J70").Value = 5777.2856
$ws.Range("L70").Value = 5777.2856
$ws.Range("N70").Value = -6317.2856

$ws.Range("H73").Value = 5425.552
$ws.Range("J73").Value = 5777.2856
$ws.Range("L73").Value = 5777.2856
$ws.Range("N73").Value = -7649.2856

$ws.Range("H102").Value = 1753
$ws.Range("I102").Value = 1562
$ws.Range("J102").Value = 5000
$ws.Range("K102").Value = 1562
$ws.Range("L102").Value = 5000
$ws.Range("M102").Value = 60
$ws.Range("N102").Value = -8244

$ws.Range("H107").Value = 2000
$ws.Range("J107").Value = 2000
$ws.Range("L107").Value = 2000
$ws.Range("N107").Value = -5840

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1855.4445
$ws.Range("I81").Value = 1333.6666
$ws.Range("J81").Value = 2116.3333
$ws.Range("K81").Value = 2667.3332
$ws.Range("L81").Value = 4232.6666
$ws.Range("M81").Value = -1606.3332
$ws.Range("N81").Value = -6354.6666

$ws.Range("H84").Value = 1855.4445
$ws.Range("I84").Value = 1333.6666
$ws.Range("J84").Value = 2116.3333
$ws.Range("K84").Value = 13336.666
$ws.Range("L84").Value = 21163.333
$ws.Range("M84").Value = -8032.666000000001
$ws.Range("N84").Value = -31771.333
